$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Couple of retired adults"

$ws.Range("B6").Select()
